# Commit message: "Pasé datos a carpeta datos"
# (i.e. "Moved the data to the 'datos' folder" - a repo reorganisation).
# Diffing the OOXML shows no actual cell/data changes (every numeric <v>
# diff is the exact same IEEE-754 double re-serialized with a different
# number of significant digits - an artifact of whichever tool/Excel
# build wrote the file, not a content edit). The one real, user-visible
# change inside the workbook itself is that the single worksheet is
# renamed from "Datos" to the generic default name "Sheet1".

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Datos") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "Sheet1"
